$p = $ppt.ActivePresentation
$s = $p.Slides.Item(21)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$para = $tr.Paragraphs(4, 1)

$newText = "random.random() ile üretilecek 100 tane sayı içinde birbirine en yakın iki sayıyı bulan programı kodlayın."
$para.Text = $newText

$para = $tr.Paragraphs(4, 1)

$lens = @(13, 27, 12, 36, 18)
$pos = 1
foreach ($len in $lens) {
    $chunk = $para.Characters($pos, $len)
    $chunk.Text = $chunk.Text
    $pos += $len
}
